# Apply the "Parallel processing implemented E" edit to the workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Algorithm" sheet: drop the separate "evaluator" row and rename the
#    remaining "evaluate" row to "life_cycle"
# ---------------------------------------------------------------------
$wsAlgorithm = $wb.Worksheets.Item("Algorithm")
$wsAlgorithm.Activate()

# Remove the old first row (name=evaluator / mj_utilities.util_general /
# evaluate_pop_parallel) entirely - the rows below shift up.
$wsAlgorithm.Rows.Item(1).Delete()

# The row that used to hold "evaluate" (now row 2) is renamed.
$wsAlgorithm.Range("A2").Value = "life_cycle"

# Widen column B a bit (manually resized by the user)
$wsAlgorithm.Columns.Item(2).ColumnWidth = 22.7109375

$wsAlgorithm.Range("A4").Select()

# ---------------------------------------------------------------------
# 2. "Project" sheet: add the new parallel-execution configuration rows
# ---------------------------------------------------------------------
$wsProject = $wb.Worksheets.Item("Project")
$wsProject.Activate()

$wsProject.Range("A10").Value = "execution"
$wsProject.Range("B10").Value = "parallel"
$wsProject.Range("A11").Value = "parallel_delay"
$wsProject.Range("B11").Value = 0
$wsProject.Range("A12").Value = "Maximum_CPU"
$wsProject.Range("B12").Value = 80
$wsProject.Range("A13").Value = "Maximum_processes"
$wsProject.Range("B13").Value = 6

# Selection moves to the newly added B12 cell
$wsProject.Range("B12").Select()

# ---------------------------------------------------------------------
# 3. "Range Variables" sheet: tighten the STEP column from 0.01 to 0.001
#    for every variable row, and make this the active sheet/selection
# ---------------------------------------------------------------------
$wsRangeVars = $wb.Worksheets.Item("Range Variables")
$wsRangeVars.Activate()

$wsRangeVars.Range("D2:D31").Value = 0.001

# Final selection/active sheet, matching the saved workbook state
$wsRangeVars.Range("D2:D31").Select()
